$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bulk-updated from
# 45192 (2023-09-23) to 45202 (2023-10-03) for every data row (2..360).
$ws.Range("C2:C360").Value = 45202
